$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 29.32133366666666
$ws.Range("H2").Value = 87.964001
$ws.Range("I2").Value = 0.006401919837078288
$ws.Range("J2").Value = 0.006401919837078288
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 8.131233999999999
$ws.Range("N2").Value = 24.393702
$ws.Range("O2").Value = 0.02090995573015822
$ws.Range("P2").Value = 0.02090995573015823
$ws.Range("Q2").Value = 238.4186252357446
$ws.Range("R2").Value = 2145.767627121701
$ws.Range("S2").Value = 0.0001338638603813287
$ws.Range("T2").Value = 0.0001338638603813288

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 29.32133366666666
$ws.Range("H3").Value = 87.964001
$ws.Range("I3").Value = 0.006401919837078288
$ws.Range("J3").Value = 0.006401919837078288
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 243.3763986666667
$ws.Range("N3").Value = 730.1291960000001
$ws.Range("O3").Value = 0.625857000534647
$ws.Range("P3").Value = 0.6258570005346471
$ws.Range("Q3").Value = 7136.120591897022
$ws.Range("R3").Value = 64225.0853270732
$ws.Range("S3").Value = 0.004006686346897073
$ws.Range("T3").Value = 0.004006686346897074

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 29.32133366666666
$ws.Range("H4").Value = 87.964001
$ws.Range("I4").Value = 0.006401919837078288
$ws.Range("J4").Value = 0.006401919837078288
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 103.9426383333333
$ws.Range("N4").Value = 311.827915
$ws.Range("O4").Value = 0.2672947262403034
$ws.Range("P4").Value = 0.2672947262403035
$ws.Range("Q4").Value = 3047.736780765324
$ws.Range("R4").Value = 27429.63102688792
$ws.Range("S4").Value = 0.001711199410264209
$ws.Range("T4").Value = 0.001711199410264209

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 29.32133366666666
$ws.Range("H5").Value = 87.964001
$ws.Range("I5").Value = 0.006401919837078288
$ws.Range("J5").Value = 0.006401919837078288
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 33.41874933333333
$ws.Range("N5").Value = 100.256248
$ws.Range("O5").Value = 0.08593831749489127
$ws.Range("P5").Value = 0.08593831749489128
$ws.Range("Q5").Value = 979.8822999253607
$ws.Range("R5").Value = 8818.940699328248
$ws.Range("S5").Value = 0.0005501702195356764
$ws.Range("T5").Value = 0.0005501702195356766

# Row 6
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 45.524413
$ws.Range("H6").Value = 136.573239
$ws.Range("I6").Value = 0.009939644832300594
$ws.Range("J6").Value = 0.009939644832300592
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 8.131233999999999
$ws.Range("N6").Value = 24.393702
$ws.Range("O6").Value = 0.02090995573015822
$ws.Range("P6").Value = 0.02090995573015823
$ws.Range("Q6").Value = 370.169654815642
$ws.Range("R6").Value = 3331.526893340778
$ws.Range("S6").Value = 0.0002078375334169014
$ws.Range("T6").Value = 0.0002078375334169014

# Row 7
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 45.524413
$ws.Range("H7").Value = 136.573239
$ws.Range("I7").Value = 0.009939644832300594
$ws.Range("J7").Value = 0.009939644832300592
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 243.3763986666667
$ws.Range("N7").Value = 730.1291960000001
$ws.Range("O7").Value = 0.625857000534647
$ws.Range("P7").Value = 0.6258570005346471
$ws.Range("Q7").Value = 11079.56768735399
$ws.Range("R7").Value = 99716.10918618586
$ws.Range("S7").Value = 0.006220796301123354
$ws.Range("T7").Value = 0.006220796301123354

# Row 8
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 45.524413
$ws.Range("H8").Value = 136.573239
$ws.Range("I8").Value = 0.009939644832300594
$ws.Range("J8").Value = 0.009939644832300592
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 103.9426383333333
$ws.Range("N8").Value = 311.827915
$ws.Range("O8").Value = 0.2672947262403034
$ws.Range("P8").Value = 0.2672947262403035
$ws.Range("Q8").Value = 4731.927595796299
$ws.Range("R8").Value = 42587.34836216669
$ws.Range("S8").Value = 0.002656814644375634
$ws.Range("T8").Value = 0.002656814644375634

# Row 9
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 45.524413
$ws.Range("H9").Value = 136.573239
$ws.Range("I9").Value = 0.009939644832300594
$ws.Range("J9").Value = 0.009939644832300592
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 33.41874933333333
$ws.Range("N9").Value = 100.256248
$ws.Range("O9").Value = 0.08593831749489127
$ws.Range("P9").Value = 0.08593831749489128
$ws.Range("Q9").Value = 1521.368946594141
$ws.Range("R9").Value = 13692.32051934727
$ws.Range("S9").Value = 0.0008541963533847038
$ws.Range("T9").Value = 0.0008541963533847037

# Row 10
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 4438.215250666667
$ws.Range("H10").Value = 13314.645752
$ws.Range("I10").Value = 0.9690247577915309
$ws.Range("J10").Value = 0.9690247577915307
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 8.131233999999999
$ws.Range("N10").Value = 24.393702
$ws.Range("O10").Value = 0.02090995573015822
$ws.Range("P10").Value = 0.02090995573015823
$ws.Range("Q10").Value = 36088.16674553933
$ws.Range("R10").Value = 324793.5007098539
$ws.Range("S10").Value = 0.0202622647868482
$ws.Range("T10").Value = 0.02026226478684821

# Row 11
$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 4438.215250666667
$ws.Range("H11").Value = 13314.645752
$ws.Range("I11").Value = 0.9690247577915309
$ws.Range("J11").Value = 0.9690247577915307
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 243.3763986666667
$ws.Range("N11").Value = 730.1291960000001
$ws.Range("O11").Value = 0.625857000534647
$ws.Range("P11").Value = 0.6258570005346471
$ws.Range("Q11").Value = 1080156.844214731
$ws.Range("R11").Value = 9721411.597932577
$ws.Range("S11").Value = 0.6064709283552203
$ws.Range("T11").Value = 0.6064709283552203

# Row 12
$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 4438.215250666667
$ws.Range("H12").Value = 13314.645752
$ws.Range("I12").Value = 0.9690247577915309
$ws.Range("J12").Value = 0.9690247577915307
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 103.9426383333333
$ws.Range("N12").Value = 311.827915
$ws.Range("O12").Value = 0.2672947262403034
$ws.Range("P12").Value = 0.2672947262403035
$ws.Range("Q12").Value = 461319.8026455297
$ws.Range("R12").Value = 4151878.223809768
$ws.Range("S12").Value = 0.2590152073539636
$ws.Range("T12").Value = 0.2590152073539636

# Row 13
$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 4438.215250666667
$ws.Range("H13").Value = 13314.645752
$ws.Range("I13").Value = 0.9690247577915309
$ws.Range("J13").Value = 0.9690247577915307
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 33.41874933333333
$ws.Range("N13").Value = 100.256248
$ws.Range("O13").Value = 0.08593831749489127
$ws.Range("P13").Value = 0.08593831749489128
$ws.Range("Q13").Value = 148319.6029494065
$ws.Range("R13").Value = 1334876.426544659
$ws.Range("S13").Value = 0.08327635729549869
$ws.Range("T13").Value = 0.08327635729549869

# Row 14
$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 67.02347933333333
$ws.Range("H14").Value = 201.070438
$ws.Range("I14").Value = 0.01463367753909034
$ws.Range("J14").Value = 0.01463367753909034
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 8.131233999999999
$ws.Range("N14").Value = 24.393702
$ws.Range("O14").Value = 0.02090995573015822
$ws.Range("P14").Value = 0.02090995573015823
$ws.Range("Q14").Value = 544.9835939534972
$ws.Range("R14").Value = 4904.852345581476
$ws.Range("S14").Value = 0.0003059895495117897
$ws.Range("T14").Value = 0.0003059895495117898

# Row 15
$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 67.02347933333333
$ws.Range("H15").Value = 201.070438
$ws.Range("I15").Value = 0.01463367753909034
$ws.Range("J15").Value = 0.01463367753909034
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 243.3763986666667
$ws.Range("N15").Value = 730.1291960000001
$ws.Range("O15").Value = 0.625857000534647
$ws.Range("P15").Value = 0.6258570005346471
$ws.Range("Q15").Value = 16311.93302625643
$ws.Range("R15").Value = 146807.3972363079
$ws.Range("S15").Value = 0.009158589531406313
$ws.Range("T15").Value = 0.009158589531406316

# Row 16
$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 67.02347933333333
$ws.Range("H16").Value = 201.070438
$ws.Range("I16").Value = 0.01463367753909034
$ws.Range("J16").Value = 0.01463367753909034
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 103.9426383333333
$ws.Range("N16").Value = 311.827915
$ws.Range("O16").Value = 0.2672947262403034
$ws.Range("P16").Value = 0.2672947262403035
$ws.Range("Q16").Value = 6966.597272186307
$ws.Range("R16").Value = 62699.37544967677
$ws.Range("S16").Value = 0.003911504831700029
$ws.Range("T16").Value = 0.003911504831700031

# Row 17
$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 67.02347933333333
$ws.Range("H17").Value = 201.070438
$ws.Range("I17").Value = 0.01463367753909034
$ws.Range("J17").Value = 0.01463367753909034
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 33.41874933333333
$ws.Range("N17").Value = 100.256248
$ws.Range("O17").Value = 0.08593831749489127
$ws.Range("P17").Value = 0.08593831749489128
$ws.Range("Q17").Value = 2239.840855288513
$ws.Range("R17").Value = 20158.56769759662
$ws.Range("S17").Value = 0.001257593626472205
$ws.Range("T17").Value = 0.001257593626472205

